# Applies numeric updates to the "想去人数" (interest count) column F
# across the 展览, 演出 and 全部类型 sheets, per the regenerated data dump.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$ws1.Range("F3").Value = 723
$ws1.Range("F5").Value = 45
$ws1.Range("F6").Value = 2890
$ws1.Range("F7").Value = 1689
$ws1.Range("F8").Value = 1914
$ws1.Range("F11").Value = 779
$ws1.Range("F12").Value = 926
$ws1.Range("F13").Value = 184
$ws1.Range("F14").Value = 392
$ws1.Range("F19").Value = 6928
$ws1.Range("F21").Value = 1680
$ws1.Range("F22").Value = 172
$ws1.Range("F23").Value = 184
$ws1.Range("F24").Value = 159
$ws1.Range("F25").Value = 340
$ws1.Range("F27").Value = 71
$ws1.Range("F28").Value = 1113
$ws1.Range("F31").Value = 102
$ws1.Range("F33").Value = 802
$ws1.Range("F35").Value = 164
$ws1.Range("F36").Value = 2
$ws1.Range("F37").Value = 147
$ws1.Range("F38").Value = 238
$ws1.Range("F39").Value = 29
$ws1.Range("F41").Value = 244
$ws1.Range("F43").Value = 178

# --- 演出 (sheet2) ---
$ws2.Range("F6").Value = 6

# --- 全部类型 (sheet4) ---
$ws4.Range("F3").Value = 723
$ws4.Range("F8").Value = 45
$ws4.Range("F9").Value = 2890
$ws4.Range("F10").Value = 1689
$ws4.Range("F11").Value = 1914
$ws4.Range("F14").Value = 779
$ws4.Range("F16").Value = 926
$ws4.Range("F17").Value = 184
$ws4.Range("F18").Value = 392
$ws4.Range("F22").Value = 6928
$ws4.Range("F24").Value = 1680
$ws4.Range("F25").Value = 6
$ws4.Range("F26").Value = 172
$ws4.Range("F27").Value = 184
$ws4.Range("F28").Value = 159
$ws4.Range("F29").Value = 340
$ws4.Range("F31").Value = 71
$ws4.Range("F32").Value = 1113
$ws4.Range("F35").Value = 102
$ws4.Range("F36").Value = 802
$ws4.Range("F38").Value = 164
$ws4.Range("F39").Value = 2
$ws4.Range("F40").Value = 147
$ws4.Range("F41").Value = 238
$ws4.Range("F42").Value = 29
$ws4.Range("F44").Value = 244
$ws4.Range("F49").Value = 178

